$wb = $excel.ActiveWorkbook

function Clone-Sheet($srcName, $newName) {
    $src = $wb.Worksheets.Item($srcName)
    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $src.Copy($null, $last)
    $newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet.Name = $newName
    return $newSheet
}

# 1. Fix the shared text used for the "choose test" instruction on Home,
#    and the "Nama Ujian" column header text used on every test tab.
$homeWs = $wb.Worksheets.Item("Home")
$homeWs.Range("B2").Value = "Sila Pilih Tab Ujian "

foreach ($name in @("IPU","IBK","IP","PPKP")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B2").Value = "Nama Ujian"
}

# 2. Select D3 on Home tab (so it is the active cell when users land there).
$homeWs.Activate()
$homeWs.Range("D3").Select() | Out-Null

# 3. Add the three new test tabs, cloned from IPU, after PPKP.
$iso = Clone-Sheet "IPU" "ISO"
$iso.Range("B2").Value = "Nama Ujian"
$iso.Range("A1:X1").Select() | Out-Null

$ukhlp = Clone-Sheet "IPU" "UKHLP"
$ukhlp.Range("B2").Value = "Nama Ujian"
$ukhlp.Range("A1:X1").Select() | Out-Null

$ukbp = Clone-Sheet "IPU" "UKBP"
$ukbp.Range("B2").Value = "Nama Ujian"
$ukbp.Range("F13").Select() | Out-Null

# 4. Leave focus back on the Home tab and save.
$homeWs.Activate()
$wb.Save()
